$d = $word.ActiveDocument

$d.Content.Find.Execute("2024-10-26 Saturday", $true, $false, $false, $false, $false, $true, 1, $false, "2024-10-27 Sunday", 2) | Out-Null
$d.Content.Find.Execute("258÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "357÷7=", 2) | Out-Null
$d.Content.Find.Execute("475÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "305÷2=", 2) | Out-Null
$d.Content.Find.Execute("427÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "293÷2=", 2) | Out-Null
$d.Content.Find.Execute("129÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "847÷5=", 2) | Out-Null
$d.Content.Find.Execute("854÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "356÷3=", 2) | Out-Null
$d.Content.Find.Execute("329÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "645÷3=", 2) | Out-Null
$d.Content.Find.Execute("212÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "385÷7=", 2) | Out-Null
$d.Content.Find.Execute("753÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "661÷9=", 2) | Out-Null
$d.Content.Find.Execute("222÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "573÷8=", 2) | Out-Null
$d.Content.Find.Execute("127÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "760÷4=", 2) | Out-Null
$d.Content.Find.Execute("545÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "717÷9=", 2) | Out-Null
$d.Content.Find.Execute("963÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "704÷6=", 2) | Out-Null
$d.Content.Find.Execute("736÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "989÷4=", 2) | Out-Null
$d.Content.Find.Execute("666÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "329÷3=", 2) | Out-Null
$d.Content.Find.Execute("654÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "798÷6=", 2) | Out-Null
$d.Content.Find.Execute("983÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "565÷6=", 2) | Out-Null
$d.Content.Find.Execute("898÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "304÷7=", 2) | Out-Null
$d.Content.Find.Execute("449÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "142÷2=", 2) | Out-Null
$d.Content.Find.Execute("470÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "364÷4=", 2) | Out-Null
$d.Content.Find.Execute("133÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "955÷7=", 2) | Out-Null
$d.Content.Find.Execute("988÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "433÷3=", 2) | Out-Null
$d.Content.Find.Execute("315÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "855÷6=", 2) | Out-Null
$d.Content.Find.Execute("298÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "971÷8=", 2) | Out-Null
$d.Content.Find.Execute("183÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "942÷5=", 2) | Out-Null
$d.Content.Find.Execute("317÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "315÷8=", 2) | Out-Null
